# Auto-generated edit script: updates cryptos price (D) and volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.885.90"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.741.20"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'230.68"
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.5256"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").Value = "'0.2750"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").Value = "'39.43"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").Value = "'0.06142"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").Value = "1.738.97"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "'0.07098"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "'15.19"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "'0.6420"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "'4.530"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "'77.49"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "'0.9998"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'0.9998"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "25.878.80"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "'0.000006681"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "1.962.70"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "'4.300"
$ws.Range("E23").Value = "  +6.07%  "
$ws.Range("D24").Value = "'8.747"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").Value = "'5.159"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "'140.17"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "'1.518"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "'1.790"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "'102.73"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'0.08305"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "'3.725"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").Value = "'3.523"
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("D34").Value = "'0.04523"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").Value = "'2.614"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").Value = "'0.9774"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D38").Value = "'2.683"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").Value = "'0.01590"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").Value = "'0.9991"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'99.96"
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("D43").Value = "'0.3863"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "'0.7318"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").Value = "'0.05336"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("D47").Value = "'0.1126"
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("D48").Value = "'6.254"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'53.71"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("D50").Value = "'30.13"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "'7.665"
$ws.Range("E51").Value = "  +3.34%  "
